$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "release/8.0.4"
$ws.Range("B7").Value = "X"
$ws.Range("C7").Value = "X"
$ws.Range("D7").Value = "X"
$ws.Range("E7").Value = "X"
